$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(3)
$ws.Range("A2").Value = "'001-Authenticate"
Write-Host ("A2 value=[" + $ws.Range("A2").Value + "]")
Write-Host ("A2 text=[" + $ws.Range("A2").Text + "]")
